# Apply the edits described by the commit "removed serving from video projects".
$d = $word.ActiveDocument

# 1) Anomaly Detection System for Image Pipeline - bullet 1: tweak wording.
$d.Content.Find.Execute(
    "Developed a system that filtered 92% of anomalies in image data pipelines, enhancing efficiency by 21%.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed a system that filtered 92% of anomalies in image data pipelines, enhancing pipeline insight efficiency by 21%.",
    2) | Out-Null

# 2) Same section - bullet 2 ("Employed model experimentation...") is replaced by
#    the wording that used to belong to bullet 3, and the old bullet 3
#    ("Kept track of multiple experiments and maintained model versions for ease
#    of delivery of the best model.") paragraph is removed entirely.
$d.Content.Find.Execute(
    "Employed model experimentation and version tracking for iterative improvements.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kept track of multiple experiments and maintained model version tracking for iterative improvements.",
    2) | Out-Null

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Kept track of multiple experiments and maintained model versions for ease of delivery of the best model. ") {
        $p.Range.Delete()
        break
    }
}

# 3) Privacy Blur Feature section: drop the "Kept track of multiple experiments
#    and maintained model versions.." bullet, and replace the deployment bullet
#    that follows it with that same sentence (single trailing period).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Kept track of multiple experiments and maintained model versions..") {
        $p.Range.Delete()
        break
    }
}

$d.Content.Find.Execute(
    "Deployed as a service with continuous monitoring, enabling easy updates, scalability and modular maintenance.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kept track of multiple experiments and maintained model versions.",
    2) | Out-Null

# 4) Visual Embedding for Retail section: merge the "enabling easy updates..."
#    phrase into the remaining deployment bullet.
$d.Content.Find.Execute(
    "Deployed as a service with continuous monitoring, ensuring robustness in a dynamic retail environment.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployed as a service with continuous monitoring, enabling easy updates, scalability and modular maintenance, ensuring robustness in a dynamic retail environment.",
    2) | Out-Null

# 5) NSFW video content filtering: replace the deployment bullet's text.
$d.Content.Find.Execute(
    "Model deployed as a service ensuring continuous delivery, easy updates, and ease of scalability",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kept track of multiple experiments and maintained model versions.",
    2) | Out-Null

# 6) Video tagging section: drop the "Automated genre classification..." bullet,
#    and move its text onto the bullet that used to describe deployment.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Automated genre classification of videos , achieving 85% efficiency and enhancing annotation processes by 90%. ") {
        $p.Range.Delete()
        break
    }
}

$d.Content.Find.Execute(
    "System deployed as a service, ensuring continuous delivery and ease of scalability",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Automated genre classification of videos , achieving 85% efficiency and enhancing annotation processes by 90%. ",
    2) | Out-Null
